$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C45").Value = "4ecce020"
$ws.Range("B45").Value = "44"
$ws.Range("D45").Value = "PascalWesterhof"
$ws.Range("E45").Value = "11-10-2024, 13.34"
$ws.Range("B46").Value = "45"
$ws.Range("F45").Value = "Basis aangelegd voor blog pagina en header hieraan toegevoegd."
$ws.Range("F46").Value = "Heb een logboek aangemaakt voor onze logs, en heb mijn overons pagina bijna afgemaakt. Ook heb ik wat extra informatie bij de style guide en de README."
$ws.Range("C46").Value = "59b8e48f"
$ws.Range("E46").Value = "11-10-2024, 15.56"
$ws.Range("D46").Value = "DinandRengers"

$ws.Range("E44").Copy()
$ws.Range("E45:E46").PasteSpecial(-4122)

[void]$ws.Range("F55:F56").Select()
